# Append scrape-run: 2025-11-10 06:29 JST
# The "ランサーズ" (Lancers) sheet is an append-only scrape log: the
# oldest row (old row 2) drops off, the surviving rows shift up by one,
# and new rows are appended at the bottom - net effect here is the sheet
# shrinks from 15 data rows (rows 2-16) down to 7 data rows (rows 2-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Drop the now-unused tail rows (old rows 9-16) -------------------
$ws.Range("A9:H16").EntireRow.Delete()

# --- Clear every existing hyperlink; we rebuild F2:F8 from scratch ---
$ws.Range("A1").Hyperlinks.Delete()

$timestamp = "2025-11-10 06:29:07"

# --- Row 2 -------------------------------------------------------------
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【募集】習慣化+目標管理を目的としたAIネイティブなWebサービスのMVP開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5430365"
$ws.Range("G2").Value = 378
$ws.Range("H2").Value = "🔥AI,Ai ◆開発 ◇管理"

# --- Row 3 -------------------------------------------------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "通話機能の安定化・不具合調査/改修(React Native × Node.js)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5430799"
$ws.Range("G3").Value = 163
$ws.Range("H3").Value = "🔥React ◆Node.js"

# --- Row 4 -------------------------------------------------------------
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "~ 5,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5416665"
$ws.Range("G4").Value = 70
$ws.Range("H4").Value = "◆効率化"

# --- Row 5 -------------------------------------------------------------
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【業務効率化】SlackとHubSpotの活用支援をお願いします"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5430436"
$ws.Range("G5").Value = 70
$ws.Range("H5").Value = "◆効率化"

# --- Row 6 (no skill-summary this time) --------------------------------
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5428756"
$ws.Range("G6").Value = 25
$ws.Range("H6").ClearContents()

# --- Row 7 (no skill-summary this time) --------------------------------
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5428755"
$ws.Range("G7").Value = 25
$ws.Range("H7").ClearContents()

# --- Row 8 (no skill-summary this time) --------------------------------
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【データ加工のプロ募集】施設情報データの修正・整備依頼"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5417622"
$ws.Range("G8").Value = 10
$ws.Range("H8").ClearContents()

# --- Rebuild the F2:F8 hyperlinks (rId order matches row order) -------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5430365")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5430799")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5416665")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5430436")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5428756")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5428755")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5417622")

# --- Column-width tweaks (D: 32 -> 30, H: 21 -> 17 "characters") -------
# ColumnWidth is in "characters"; the stored xml width adds ~5/6 of a
# character of internal padding, so back that out to land exactly on
# the target xml width.
$ws.Columns.Item(4).ColumnWidth = 29.16666666666667
$ws.Columns.Item(8).ColumnWidth = 16.16666666666667
